$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 97; $row++) {
    $newIndex = $row - 2
    $ws.Cells.Item($row, 1).Value = "q$newIndex"
}
